$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").Value = ""

$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("M14").Value = ""

$ws.Range("H53").Value = 240.15384
$ws.Range("I53").Value = 21
$ws.Range("K53").Value = 21
$ws.Range("M53").Value = 616

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 446465.2
$ws.Range("I32").Value = 461700.2
$ws.Range("J32").Value = 334742
$ws.Range("K32").Value = 461700.2
$ws.Range("L32").Value = 334742
$ws.Range("M32").Value = -461413.2
$ws.Range("N32").Value = -335316

$ws.Range("H45").Value = 2008.6086
$ws.Range("I45").Value = 1791
$ws.Range("K45").Value = 1791
$ws.Range("M45").Value = -1414

$ws.Range("H61").Value = 2990.2942
$ws.Range("I61").Value = 1802.5
$ws.Range("J61").Value = 8533.333000000001
$ws.Range("K61").Value = 1802.5
$ws.Range("L61").Value = 8533.333000000001
$ws.Range("M61").Value = -1590.5
$ws.Range("N61").Value = -8957.333000000001

$ws.Range("H74").Value = 49804.27
$ws.Range("I74").Value = 54983.11
$ws.Range("J74").Value = 1900
$ws.Range("K74").Value = 54983.11
$ws.Range("L74").Value = 1900
$ws.Range("M74").Value = -54109.11
$ws.Range("N74").Value = -3648

$ws.Range("H77").Value = 49804.27
$ws.Range("I77").Value = 54983.11
$ws.Range("J77").Value = 1900
$ws.Range("K77").Value = 274915.55
$ws.Range("L77").Value = 9500
$ws.Range("M77").Value = -270547.55
$ws.Range("N77").Value = -18236

$ws.Range("H122").Value = 1458.2
$ws.Range("I122").Value = 1397.4286
$ws.Range("K122").Value = 4192.2858
$ws.Range("M122").Value = -1742.2858

$ws.Range("H132").Value = 3867722
$ws.Range("I132").Value = 4638224
$ws.Range("J132").Value = 1446143.8
$ws.Range("K132").Value = 13914672
$ws.Range("L132").Value = 4338431.4
$ws.Range("M132").Value = -13912142
$ws.Range("N132").Value = -4343491.4

$ws.Range("H136").Value = 2990.2942
$ws.Range("I136").Value = 1802.5
$ws.Range("J136").Value = 8533.333000000001
$ws.Range("K136").Value = 5407.5
$ws.Range("L136").Value = 25599.999
$ws.Range("M136").Value = -2857.5
$ws.Range("N136").Value = -30699.999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 13928.571
$ws.Range("I26").Value = 9375
$ws.Range("J26").Value = 20000
$ws.Range("K26").Value = 9375
$ws.Range("L26").Value = 20000
$ws.Range("M26").Value = -9083
$ws.Range("N26").Value = -20584

$ws.Range("H28").Value = 40000
$ws.Range("J28").Value = 40000
$ws.Range("L28").Value = 40000
$ws.Range("N28").Value = -40588

$ws.Range("H57").Value = 36499.5
$ws.Range("I57").Value = 1000
$ws.Range("J57").Value = 71999
$ws.Range("K57").Value = 1000
$ws.Range("L57").Value = 71999
$ws.Range("M57").Value = -280
$ws.Range("N57").Value = -73439

$ws.Range("H60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").Value = ""

$ws.Range("H105").Value = 0
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = ""
$ws.Range("N105").Value = ""

$ws.Range("H136").Value = 36499.5
$ws.Range("I136").Value = 1000
$ws.Range("J136").Value = 71999
$ws.Range("K136").Value = 1000
$ws.Range("L136").Value = 71999
$ws.Range("M136").Value = 4100
$ws.Range("N136").Value = -82199

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").Value = ""

$ws.Range("H103").Value = 9108.777
$ws.Range("I103").Value = 5666.5
$ws.Range("J103").Value = 15993.333
$ws.Range("K103").Value = 5666.5
$ws.Range("L103").Value = 15993.333
$ws.Range("M103").Value = -4494.5
$ws.Range("N103").Value = -18337.333

$ws.Range("H108").Value = 17310.5
$ws.Range("I108").Value = 17310.5
$ws.Range("J108").Value = 0
$ws.Range("K108").Value = 17310.5
$ws.Range("L108").Value = 0
$ws.Range("M108").Value = -13470.5
$ws.Range("N108").Value = ""

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H16").Value = 600
$ws.Range("I16").Value = 600
$ws.Range("K16").Value = 1800
$ws.Range("M16").Value = -1627

$ws.Range("H126").Value = 5176.6665
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 5176.6665
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 15529.9995
$ws.Range("M126").Value = ""
$ws.Range("N126").Value = -25409.9995

$ws.Range("H136").Value = 1745.9474
$ws.Range("I136").Value = 994
$ws.Range("J136").Value = 2581.4443
$ws.Range("K136").Value = 2982
$ws.Range("L136").Value = 7744.3329
$ws.Range("M136").Value = 2118
$ws.Range("N136").Value = -17944.3329

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4233.3335
$ws.Range("I70").Value = 4190.476
$ws.Range("J70").Value = 4333.3335
$ws.Range("K70").Value = 4190.476
$ws.Range("L70").Value = 4333.3335
$ws.Range("M70").Value = -3920.476
$ws.Range("N70").Value = -4873.3335

$ws.Range("H73").Value = 4233.3335
$ws.Range("I73").Value = 4190.476
$ws.Range("J73").Value = 4333.3335
$ws.Range("K73").Value = 4190.476
$ws.Range("L73").Value = 4333.3335
$ws.Range("M73").Value = -3254.476
$ws.Range("N73").Value = -6205.3335

$ws.Range("H102").Value = 23593.111
$ws.Range("I102").Value = 8110.8
$ws.Range("J102").Value = 101004.664
$ws.Range("K102").Value = 8110.8
$ws.Range("L102").Value = 101004.664
$ws.Range("M102").Value = -6488.8
$ws.Range("N102").Value = -104248.664

$ws.Range("H113").Value = 1615.6111
$ws.Range("I113").Value = 1645.1428
$ws.Range("J113").Value = 1512.25
$ws.Range("K113").Value = 1645.1428
$ws.Range("L113").Value = 1512.25
$ws.Range("M113").Value = 524.8571999999999
$ws.Range("N113").Value = -5852.25

$ws.Range("H122").Value = 3856.84
$ws.Range("I122").Value = 2936.3572
$ws.Range("J122").Value = 5028.364
$ws.Range("K122").Value = 8809.071599999999
$ws.Range("L122").Value = 15085.092
$ws.Range("M122").Value = -6359.071599999999
$ws.Range("N122").Value = -19985.092

$ws.Range("H132").Value = 62897.883
$ws.Range("I132").Value = 2950
$ws.Range("J132").Value = 70890.92999999999
$ws.Range("K132").Value = 8850
$ws.Range("L132").Value = 212672.79
$ws.Range("M132").Value = -6320
$ws.Range("N132").Value = -217732.79

$ws.Range("H137").Value = 57786.668
$ws.Range("J137").Value = 57786.668
$ws.Range("L137").Value = 57786.668
$ws.Range("N137").Value = -67986.66800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 46361.523
$ws.Range("I40").Value = 1883.3334
$ws.Range("K40").Value = 1883.3334
$ws.Range("M40").Value = -1747.3334

$ws.Range("H109").Value = 13000.8
$ws.Range("J109").Value = 13000.8
$ws.Range("L109").Value = 13000.8
$ws.Range("N109").Value = -15774.8

$ws.Range("H122").Value = 2098.6155
$ws.Range("I122").Value = 1178.8572
$ws.Range("J122").Value = 3171.6667
$ws.Range("K122").Value = 3536.5716
$ws.Range("L122").Value = 9515.000100000001
$ws.Range("M122").Value = -1086.5716
$ws.Range("N122").Value = -14415.0001

$ws.Range("H132").Value = 634579.25
$ws.Range("I132").Value = 145664.58
$ws.Range("J132").Value = 2003540.4
$ws.Range("K132").Value = 436993.74
$ws.Range("L132").Value = 6010621.199999999
$ws.Range("M132").Value = -434463.74
$ws.Range("N132").Value = -6015681.199999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 961.8421
$ws.Range("I107").Value = 694.4286
$ws.Range("J107").Value = 1710.6
$ws.Range("K107").Value = 2083.2858
$ws.Range("L107").Value = 5131.799999999999
$ws.Range("M107").Value = -163.2857999999997
$ws.Range("N107").Value = -8971.799999999999

$ws.Range("H132").Value = 4997.3335
$ws.Range("I132").Value = 1500.7222
$ws.Range("J132").Value = 10242.25
$ws.Range("K132").Value = 4502.1666
$ws.Range("L132").Value = 30726.75
$ws.Range("M132").Value = -1972.1666
$ws.Range("N132").Value = -35786.75
